$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the last row's FirstName/LastName values
$ws.Range("A4").Value = "Zman"
$ws.Range("B4").Value = "Zach"

# Update the active selection to match the edited cell
$ws.Range("B4").Select()
